# Auto-generated edit script: refresh crypto price/volume columns (D, E)
# matching the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.208.10'
$ws.Range("E2").Value = '  -1.97%  '

$ws.Range("D3").Value = '1.574.11'
$ws.Range("E3").Value = '  -1.29%  '

$ws.Range("E4").Value = '  -0.39%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '207.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.14%  '

$ws.Range("E6").Value = '  -2.31%  '

$ws.Range("E7").Value = '  -0.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.36'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("E9").Value = '  -1.93%  '

$ws.Range("E11").Value = '  -0.50%  '

$ws.Range("D12").Value = '1.796.66'
$ws.Range("E12").Value = '  -1.34%  '

$ws.Range("D13").Value = '1.574.50'
$ws.Range("E13").Value = '  -1.33%  '

$ws.Range("E14").Value = '  -2.04%  '

$ws.Range("E15").Value = '  -2.33%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.24%  '

$ws.Range("D17").Value = '27.194.83'
$ws.Range("E17").Value = '  -2.02%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '215.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.92%  '

$ws.Range("D20").Value = '0.0₃0686'
$ws.Range("E20").Value = '  -1.44%  '

$ws.Range("E21").Value = '  -0.34%  '

$ws.Range("E22").Value = '  -0.44%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.45%  '

$ws.Range("E24").Value = '  +0.70%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.76'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.62%  '

$ws.Range("E26").Value = '  -6.56%  '

$ws.Range("E27").Value = '  -1.44%  '

$ws.Range("E28").Value = '  -1.26%  '

$ws.Range("E29").Value = '  -0.33%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.70%  '

$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("E32").Value = '  -1.73%  '

$ws.Range("D33").Value = '1.400.78'
$ws.Range("E33").Value = '  +1.38%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.77%  '

$ws.Range("E35").Value = '  +1.18%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.945'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.95%  '

$ws.Range("E37").Value = '  -2.28%  '

$ws.Range("E38").Value = '  -1.75%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.818'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.14%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.519'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.27%  '

$ws.Range("E41").Value = '  -0.31%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.990'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.29%  '

$ws.Range("E43").Value = '  +3.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.35'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.18%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '63.82'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.15%  '

$ws.Range("E46").Value = '  +0.10%  '

$ws.Range("D47").Value = '1.709.68'
$ws.Range("E47").Value = '  -1.29%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.00'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.02%  '

$ws.Range("D49").Value = '0.0₇0989'
$ws.Range("E49").Value = '  -2.43%  '

$ws.Range("E50").Value = '  -1.29%  '

$ws.Range("E51").Value = '  -0.37%  '

